$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column M (26-jun) with a header matching the L column style ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Range("L1").Copy()
$ws1.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("M1").Value = "26-jun"

# Daily spot prices for 26-jun (column M, rows 2-25)
$ws1.Cells.Item(2, 13).Value = 119.95
$ws1.Cells.Item(3, 13).Value = 99.63
$ws1.Cells.Item(4, 13).Value = 93.02
$ws1.Cells.Item(5, 13).Value = 76.2
$ws1.Cells.Item(6, 13).Value = 77.45
$ws1.Cells.Item(7, 13).Value = 79.62
$ws1.Cells.Item(8, 13).Value = 84.95
$ws1.Cells.Item(9, 13).Value = 114.4
$ws1.Cells.Item(10, 13).Value = 106.79
$ws1.Cells.Item(11, 13).Value = 85
$ws1.Cells.Item(12, 13).Value = 69.02
$ws1.Cells.Item(13, 13).Value = 64.44
$ws1.Cells.Item(14, 13).Value = 62.59
$ws1.Cells.Item(15, 13).Value = 45.64
$ws1.Cells.Item(16, 13).Value = 37.39
$ws1.Cells.Item(17, 13).Value = 36.2
$ws1.Cells.Item(18, 13).Value = 37.39
$ws1.Cells.Item(19, 13).Value = 59.13
$ws1.Cells.Item(20, 13).Value = 80.09
$ws1.Cells.Item(21, 13).Value = 97.81
$ws1.Cells.Item(22, 13).Value = 114.53
$ws1.Cells.Item(23, 13).Value = 110.61
$ws1.Cells.Item(24, 13).Value = 112.12
$ws1.Cells.Item(25, 13).Value = 99.84

# --- Sheet "Gaz": append row 8 (2025-06-25, 35.05) ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A8").Value = "'2025-06-25"  # leading apostrophe forces text, matching existing date-as-text cells
$ws2.Range("A8").Style = "Normal"  # drop the quote-prefix formatting so the cell matches its siblings
$ws2.Range("B8").Value = 35.05

# --- Sheet "CO2": append row 8 (2025-06-25, 70.17) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A8").Value = "'2025-06-25"
$ws3.Range("A8").Style = "Normal"
$ws3.Range("B8").Value = 70.17

